$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(256.56, 1, 74.8, "2016-06-06 12:18:37", 1),
    @(256.56, 1, 74.8, "2016-06-06 14:22:05", 1),
    @(256.56, 1, 74.8, "2016-06-06 14:22:15", 1),
    @(256.56, 1, 74.8, "2016-06-06 14:22:25", 1),
    @(256.56, 1, 74.8, "2016-06-06 14:25:16", 1),
    @(256.56, 0, 77,   "2016-06-06 14:36:57", 1)
)

$r = 16
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
